{"js": "// Commit: \"changed 2 files in main\"\n// Target: append \" (Changed main)\" to the end of the paragraph that\n// reads \"This is a Microsoft word document.\", as three separate runs:\n//   \" (\", \"Changed main\", \")\"\n// (matching the literal OOXML of the target diff exactly). A plain\n// paragraph.insertText()/range.insertText() call would merge the new\n// text into the existing run (same formatting => coalesced run), so\n// we splice exact OOXML runs in with insertOoxml(), which preserves\n// run boundaries verbatim.\n\n// Locate the target paragraph by its known text so the script is not\n// dependent on a hard-coded paragraph index.\nconst searchResults = context.document.body.search(\n  \"This is a Microsoft word document.\",\n  { matchCase: true, matchWholeWord: false }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find the target paragraph (\"This is a Microsoft word document.\").');\n}\n\nconst targetParagraph = searchResults.items[0].paragraphs.getFirst();\n\nconst ooxml = [\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>',\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">',\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>',\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">',\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>',\n  '</Relationships>',\n  '</pkg:xmlData></pkg:part>',\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>',\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">',\n  '<w:body><w:p>',\n  '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>',\n  '<w:r><w:t>Changed main</w:t></w:r>',\n  '<w:r><w:t>)</w:t></w:r>',\n  '</w:p></w:body>',\n  '</w:document>',\n  '</pkg:xmlData></pkg:part>',\n  '</pkg:package>'\n].join('');\n\ntargetParagraph.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Commit: \"changed 2 files in main\"\n# Target: append \" (Changed main)\" to the end of the paragraph that\n# reads \"This is a Microsoft word document.\", as three separate runs:\n#   \" (\", \"Changed main\", \")\"\n# (matching the literal OOXML of the target diff exactly). Plain\n# Range.InsertAfter() text insertion merges the new text into the\n# existing run (identical formatting => the runs get coalesced), so\n# instead we read the paragraph's own WordOpenXML, splice in the three\n# literal <w:r> elements before its closing </w:p>, and feed that whole\n# paragraph back in via Range.InsertXML() on the paragraph's full\n# range. InsertXML() replaces the range with exactly the XML supplied,\n# so the run boundaries come out verbatim and the paragraph's original\n# attributes (w14:paraId, rsids, ...) are preserved because we read\n# them straight out of the document instead of hard-coding them.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the target sentence (rather than\n# assuming a fixed paragraph index).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"This is a Microsoft word document.*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw 'Could not find the target paragraph (\"This is a Microsoft word document.\").'\n}\n\n$r = $target.Range\n$owx = $r.WordOpenXML\n\nif ($owx -match '(<w:p[ >][\\s\\S]*?</w:p>)') {\n    $paraXml = $matches[1]\n} else {\n    throw \"Could not extract paragraph XML from WordOpenXML.\"\n}\n\n$newRuns = '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n           '<w:r><w:t>Changed main</w:t></w:r>' +\n           '<w:r><w:t>)</w:t></w:r>'\n$updatedParaXml = $paraXml -replace '</w:p>$', ($newRuns + '</w:p>')\n\n$pkg = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships></pkg:xmlData></pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  '<w:body>' + $updatedParaXml + '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($pkg)\n"}
